$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '28.026.33'
$ws.Range("E2").Value = '  -0.66%  '

# Row 3
$ws.Range("D3").Value = '1.798.87'
$ws.Range("E3").Value = '  -0.13%  '

# Row 4
$ws.Range("E4").Value = '  +0.02%  '

# Row 5
$ws.Range("D5").Value = '309.87'
$ws.Range("E5").Value = '  -1.42%  '

# Row 6
$ws.Range("D6").Value = '1.003'
$ws.Range("E6").Value = '  +0.05%  '

# Row 7
$ws.Range("D7").Value = '0.5072'
$ws.Range("E7").Value = '  -3.61%  '

# Row 8
$ws.Range("D8").Value = '0.3841'
$ws.Range("E8").Value = '  +0.53%  '

# Row 9
$ws.Range("D9").Value = '0.07694'
$ws.Range("E9").Value = '  -3.49%  '

# Row 10
$ws.Range("D10").Value = '1.093'
$ws.Range("E10").Value = '  -0.42%  '

# Row 11
$ws.Range("D11").Value = '40.72'
$ws.Range("E11").Value = '  -1.63%  '

# Row 12
$ws.Range("D12").Value = '6.311'
$ws.Range("E12").Value = '  -0.01%  '

# Row 13
$ws.Range("E13").Value = '  +0.01%  '

# Row 14
$ws.Range("D14").Value = '20.26'
$ws.Range("E14").Value = '  -1.84%  '

# Row 15
$ws.Range("D15").Value = '1.804.91'
$ws.Range("E15").Value = '  +0.01%  '

# Row 16
$ws.Range("D16").Value = '7.245'
$ws.Range("E16").Value = '  -0.88%  '

# Row 17
$ws.Range("D17").Value = '91.96'
$ws.Range("E17").Value = '  -0.96%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001068'
$ws.Range("E18").Value = '  -2.49%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06570'
$ws.Range("E19").Value = '  -0.52%  '

# Row 20
$ws.Range("D20").Value = '1.003'
$ws.Range("E20").Value = '  +0.04%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.20'
$ws.Range("E21").Value = '  -0.76%  '

# Row 22
$ws.Range("D22").Value = '5.958'
$ws.Range("E22").Value = '  -0.12%  '

# Row 23
$ws.Range("D23").Value = '28.036.99'
$ws.Range("E23").Value = '  -0.75%  '

# Row 24
$ws.Range("D24").Value = '11.02'
$ws.Range("E24").Value = '  -1.56%  '

# Row 25
$ws.Range("D25").Value = '2.233'
$ws.Range("E25").Value = '  -0.11%  '

# Row 26
$ws.Range("D26").Value = '159.92'
$ws.Range("E26").Value = '  +0.33%  '

# Row 27
$ws.Range("B27").Value = 'LidoDAOToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D27").Value = '2.413'
$ws.Range("E27").Value = '  +1.06%  '

# Row 28
$ws.Range("B28").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C28").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D28").Value = '2.011.56'
$ws.Range("E28").Value = '  +0.23%  '

# Row 29
$ws.Range("D29").Value = '20.21'
$ws.Range("E29").Value = '  -1.39%  '

# Row 30
$ws.Range("D30").Value = '127.07'
$ws.Range("E30").Value = '  +3.19%  '

# Row 31
$ws.Range("D31").Value = '0.1088'
$ws.Range("E31").Value = '  -0.59%  '

# Row 32
$ws.Range("D32").Value = '1.042'
$ws.Range("E32").Value = '  -1.58%  '

# Row 33
$ws.Range("D33").Value = '3.647'
$ws.Range("E33").Value = '  -0.37%  '

# Row 34
$ws.Range("D34").Value = '5.531'
$ws.Range("E34").Value = '  -0.03%  '

# Row 35
$ws.Range("D35").Value = '0.06961'
$ws.Range("E35").Value = '  -4.55%  '

# Row 36
$ws.Range("D36").Value = '9.073'
$ws.Range("E36").Value = '  +2.45%  '

# Row 37
$ws.Range("E37").Value = '  +1.12%  '

# Row 38
$ws.Range("D38").Value = '0.2157'
$ws.Range("E38").Value = '  -0.16%  '

# Row 39
$ws.Range("B39").Value = 'Aptos'
$ws.Range("C39").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D39").Value = '11.47'
$ws.Range("E39").Value = '  -6.11%  '

# Row 40
$ws.Range("B40").Value = 'InternetComputer(DFINITY)'
$ws.Range("C40").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D40").Value = '4.989'
$ws.Range("E40").Value = '  -1.54%  '

# Row 41
$ws.Range("D41").Value = '0.6094'
$ws.Range("E41").Value = '  -1.65%  '

# Row 42
$ws.Range("E42").Value = '  +0.02%  '

# Row 43
$ws.Range("D43").Value = '1.147'
$ws.Range("E43").Value = '  -1.49%  '

# Row 44
$ws.Range("D44").Value = '13.13'
$ws.Range("E44").Value = '  -0.72%  '

# Row 45
$ws.Range("D45").Value = '1.297'
$ws.Range("E45").Value = '  -5.23%  '

# Row 46
$ws.Range("D46").Value = '3.706'
$ws.Range("E46").Value = '  -1.39%  '

# Row 47
$ws.Range("D47").Value = '0.5869'
$ws.Range("E47").Value = '  -2.02%  '

# Row 48
$ws.Range("D48").Value = '125.22'
$ws.Range("E48").Value = '  -1.27%  '

# Row 49
$ws.Range("D49").Value = '1.185'
$ws.Range("E49").Value = '  -1.73%  '

# Row 50
$ws.Range("D50").Value = '1.918'
$ws.Range("E50").Value = '  -0.13%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06720'
$ws.Range("E51").Value = '  -1.56%  '
